# 2017-10-12_trik_cohort_list_free_run.xlsx
# modwt analysis and trikinetics free run cohort list update
#
# Extends the Trik_cohort_date / Trik_cohort_day / Day8 table on Sheet1
# from row 37 (2017-10-31) through row 67 (2017-11-30), continuing the
# same day-count and +8-day offset pattern already present in the sheet,
# and updates the saved view state (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(38, "2017-11-01", "37", "2017-11-09"),
    @(39, "2017-11-02", "38", "2017-11-10"),
    @(40, "2017-11-03", "39", "2017-11-11"),
    @(41, "2017-11-04", "40", "2017-11-12"),
    @(42, "2017-11-05", "41", "2017-11-13"),
    @(43, "2017-11-06", "42", "2017-11-14"),
    @(44, "2017-11-07", "43", "2017-11-15"),
    @(45, "2017-11-08", "44", "2017-11-16"),
    @(46, "2017-11-09", "45", "2017-11-17"),
    @(47, "2017-11-10", "46", "2017-11-18"),
    @(48, "2017-11-11", "47", "2017-11-19"),
    @(49, "2017-11-12", "48", "2017-11-20"),
    @(50, "2017-11-13", "49", "2017-11-21"),
    @(51, "2017-11-14", "50", "2017-11-22"),
    @(52, "2017-11-15", "51", "2017-11-23"),
    @(53, "2017-11-16", "52", "2017-11-24"),
    @(54, "2017-11-17", "53", "2017-11-25"),
    @(55, "2017-11-18", "54", "2017-11-26"),
    @(56, "2017-11-19", "55", "2017-11-27"),
    @(57, "2017-11-20", "56", "2017-11-28"),
    @(58, "2017-11-21", "57", "2017-11-29"),
    @(59, "2017-11-22", "58", "2017-11-30"),
    @(60, "2017-11-23", "59", "2017-12-01"),
    @(61, "2017-11-24", "60", "2017-12-02"),
    @(62, "2017-11-25", "61", "2017-12-03"),
    @(63, "2017-11-26", "62", "2017-12-04"),
    @(64, "2017-11-27", "63", "2017-12-05"),
    @(65, "2017-11-28", "64", "2017-12-06"),
    @(66, "2017-11-29", "65", "2017-12-07"),
    @(67, "2017-11-30", "66", "2017-12-08")
)

# Fill column-by-column (A down, then B down, then C down) to mirror how
# the author extended the table (matches the shared-string append order).
foreach ($row in $rows) {
    $ws.Cells.Item($row[0], 1).Value = $row[1]
}
foreach ($row in $rows) {
    $ws.Cells.Item($row[0], 2).Value = $row[2]
}
foreach ($row in $rows) {
    $ws.Cells.Item($row[0], 3).Value = $row[3]
}

# Scroll the saved view down to the newly added rows and move the
# active selection, matching the author's final view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
$ws.Range("M69").Select()

# Reposition the Excel window itself (best-effort; mirrors the author's
# final on-screen window position captured in workbookView).
$win.Left = 23220
$win.Top = 11020
